$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The account-statement data table previously held 15 rows (worker x period),
# including a period "1606" row for JOSE CARLO LEON BONILLA that is being
# removed in this update. Delete that row, shifting everything below it up.
$ws.Rows("22:22").Delete()

# New, re-sorted (by period, then by worker) data for the remaining 14 rows.
$data = @(
    @("CC", "73088156",   "OSVALDO LEON PEROZA",     "1604", 27600, 781242),
    @("CC", "1047475488", "JOSE CARLO LEON BONILLA", "1604", 27600, 689455),
    @("CC", "73102108",   "GERMAN RAMOS MANJARRES",  "1604", 27600, 781242),
    @("CC", "73088156",   "OSVALDO LEON PEROZA",     "1605", 27600, 781242),
    @("CC", "1047475488", "JOSE CARLO LEON BONILLA", "1605", 27600, 689455),
    @("CC", "73102108",   "GERMAN RAMOS MANJARRES",  "1605", 27600, 781242),
    @("CC", "73088156",   "OSVALDO LEON PEROZA",     "1606", 27600, 781242),
    @("CC", "73102108",   "GERMAN RAMOS MANJARRES",  "1606", 27600, 781242),
    @("CC", "73088156",   "OSVALDO LEON PEROZA",     "1607", 27578, 781242),
    @("CC", "73102108",   "GERMAN RAMOS MANJARRES",  "1607", 27578, 781242),
    @("CC", "73088156",   "OSVALDO LEON PEROZA",     "1608", 27578, 781242),
    @("CC", "73102108",   "GERMAN RAMOS MANJARRES",  "1608", 27578, 781242),
    @("CC", "73088156",   "OSVALDO LEON PEROZA",     "1609", 27578, 781242),
    @("CC", "73102108",   "GERMAN RAMOS MANJARRES",  "1609", 27578, 781242)
)

$row = 16
$totalMora = 0
foreach ($r in $data) {
    $ws.Cells.Item($row, 2).Value = $r[0]
    $ws.Cells.Item($row, 3).Value = $r[1]
    $ws.Cells.Item($row, 4).Value = $r[2]
    $ws.Cells.Item($row, 5).Value = $r[3]
    $ws.Cells.Item($row, 6).Value = $r[4]
    $ws.Cells.Item($row, 7).Value = $r[5]
    $totalMora = $totalMora + $r[4]
    $row = $row + 1
}

# Update the total "VALOR MORA" figure to reflect the refreshed data set.
$ws.Range("E11").Value = $totalMora

# Keep the bestFit column widths in sync with the refreshed content.
$ws.Columns("B:J").AutoFit()
